$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value2 = $val
    $c.Style = "Normal"
}

# Row 2
Set-TextCell 2 4 '64.263.41'
Set-TextCell 2 5 '  -2.62%  '

# Row 3
Set-TextCell 3 4 '3.332.34'
Set-TextCell 3 5 '  -4.42%  '

# Row 5
Set-TextCell 5 4 '551.29'
Set-TextCell 5 5 '  -5.43%  '

# Row 6
Set-TextCell 6 4 '174.76'
Set-TextCell 6 5 '  -1.77%  '

# Row 7
Set-TextCell 7 5 '  -2.59%  '

# Row 8
Set-TextCell 8 2 'USDC'
Set-TextCell 8 3 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
Set-TextCell 8 4 '1.00'
Set-TextCell 8 5 '  +0.06%  '

# Row 9
Set-TextCell 9 2 'LidoStakedEther'
Set-TextCell 9 3 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth'
Set-TextCell 9 4 '3.324.50'
Set-TextCell 9 5 '  -4.40%  '

# Row 10
Set-TextCell 10 2 'Cardano'
Set-TextCell 10 3 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada'
Set-TextCell 10 4 '0.625'
Set-TextCell 10 5 '  -1.59%  '

# Row 11
Set-TextCell 11 2 'Dogecoin'
Set-TextCell 11 3 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
Set-TextCell 11 4 '0.163'
Set-TextCell 11 5 '  +1.01%  '

# Row 12
Set-TextCell 12 4 '54.29'
Set-TextCell 12 5 '  -2.97%  '

# Row 13
Set-TextCell 13 4 '0.0000271'
Set-TextCell 13 5 '  -2.16%  '

# Row 14
Set-TextCell 14 4 '9.03'
Set-TextCell 14 5 '  -2.46%  '

# Row 15
Set-TextCell 15 4 '3.879.71'
Set-TextCell 15 5 '  -4.00%  '

# Row 16
Set-TextCell 16 4 '18.25'
Set-TextCell 16 5 '  -0.16%  '

# Row 17
Set-TextCell 17 4 '0.117'
Set-TextCell 17 5 '  -3.16%  '

# Row 18
Set-TextCell 18 4 '3.342.53'
Set-TextCell 18 5 '  -4.42%  '

# Row 19
Set-TextCell 19 4 '64.290.99'
Set-TextCell 19 5 '  -2.49%  '

# Row 20
Set-TextCell 20 4 '11.70'
Set-TextCell 20 5 '  -3.03%  '

# Row 21
Set-TextCell 21 4 '0.976'
Set-TextCell 21 5 '  -3.39%  '

# Row 22
Set-TextCell 22 4 '431.67'
Set-TextCell 22 5 '  +4.62%  '

# Row 23
Set-TextCell 23 4 '5.10'
Set-TextCell 23 5 '  +16.66%  '

# Row 24
Set-TextCell 24 4 '4.05'
Set-TextCell 24 5 '  -5.78%  '

# Row 25
Set-TextCell 25 4 '84.14'
Set-TextCell 25 5 '  -0.71%  '

# Row 26
Set-TextCell 26 4 '13.37'
Set-TextCell 26 5 '  -0.42%  '

# Row 27
Set-TextCell 27 4 '10.71'
Set-TextCell 27 5 '  -3.15%  '

# Row 28
Set-TextCell 28 4 '2.81'
Set-TextCell 28 5 '  -1.69%  '

# Row 29
Set-TextCell 29 4 '8.71'
Set-TextCell 29 5 '  -5.33%  '

# Row 30
Set-TextCell 30 4 '29.66'
Set-TextCell 30 5 '  -1.90%  '

# Row 31
Set-TextCell 31 5 '  -0.64%  '

# Row 32
Set-TextCell 32 4 '11.46'
Set-TextCell 32 5 '  -2.35%  '

# Row 33
Set-TextCell 33 4 '580.21'
Set-TextCell 33 5 '  -2.15%  '

# Row 34
Set-TextCell 34 5 '  -3.14%  '

# Row 35
Set-TextCell 35 4 '58.18'
Set-TextCell 35 5 '  -4.51%  '

# Row 36
Set-TextCell 36 4 '1.00'
Set-TextCell 36 5 '  -0.01%  '

# Row 37
Set-TextCell 37 4 '0.141'
Set-TextCell 37 5 '  -7.88%  '

# Row 38
Set-TextCell 38 4 '3.47'
Set-TextCell 38 5 '  -3.41%  '

# Row 39
Set-TextCell 39 4 '35.53'
Set-TextCell 39 5 '  -3.67%  '

# Row 40
Set-TextCell 40 4 '0.0₃0748'
Set-TextCell 40 5 '  -5.87%  '

# Row 41
Set-TextCell 41 4 '0.366'
Set-TextCell 41 5 '  -4.76%  '

# Row 42
Set-TextCell 42 4 '3.099.73'
Set-TextCell 42 5 '  -3.83%  '

# Row 43
Set-TextCell 43 5 '  +0.22%  '

# Row 44
Set-TextCell 44 4 '2.79'
Set-TextCell 44 5 '  -5.89%  '

# Row 45
Set-TextCell 45 4 '3.22'
Set-TextCell 45 5 '  -2.33%  '

# Row 46
Set-TextCell 46 4 '0.0407'
Set-TextCell 46 5 '  -2.70%  '

# Row 47
Set-TextCell 47 4 '2.45'
Set-TextCell 47 5 '  -3.90%  '

# Row 48
Set-TextCell 48 5 '  -2.34%  '

# Row 49
Set-TextCell 49 5 '  -2.28%  '

# Row 50
Set-TextCell 50 4 '134.63'
Set-TextCell 50 5 '  -3.64%  '

# Row 51
Set-TextCell 51 4 '8.23'
Set-TextCell 51 5 '  -4.28%  '
